# Regionalization.xlsx update: add USA regional data and switch active sheet/selection
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # CAN
$ws2 = $wb.Worksheets.Item(2)   # USA

# Populate the USA sheet with region codes (column A)
$usaData = @("REGION","NW","CA","MN","SW","CE","TX","MW","AL","MA","SE","FL","NY","NE")
for ($i = 0; $i -lt $usaData.Length; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 1).Value = $usaData[$i]
}

# Update selection on CAN sheet (no longer the active tab)
$ws1.Range("D36").Select()

# Make USA the active/selected sheet, with its own selection
$ws2.Activate()
$ws2.Range("J7").Select()
